$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, using the same header style as the
# other header cells (e.g. G1: bold, bordered, centered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the corresponding value for row 2 in the new "Save" column.
$ws.Range("H2").Value = 1
